$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (pushes the former rows 52-53 down to 53-54),
# inheriting formatting from the row above as Excel normally does.
$ws.Rows("52").Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 44509
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 100112026
$ws.Range("G52").Value = "Haba"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 600
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = 6000
$ws.Range("N52").Value = "$/saco 25 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 240
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
